$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -4736178660.47286
$ws.Range("C2").Value = -4743118513.867647

$ws.Range("B3").Value = -4726374200.22957
$ws.Range("C3").Value = -4740253936.489744

$ws.Range("B4").Value = -4706765166.463247
$ws.Range("C4").Value = -4734524756.866552

$ws.Range("B5").Value = -4683234126.569675
$ws.Range("C5").Value = -4727649697.551818
